# Append BTC/USDC trade rows 30-38 (2025-06-10 session) to the Trades sheet,
# and remove the now-obsolete OBV-based rows that had been filtered out upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trades")

# Date columns (I = Entry Date, K = Exit Date) are stored as plain text "YYYY-MM-DD"
# in this report, not as Excel date serials. Force text format before assignment so
# COM does not auto-coerce the string into a date value, then restore Normal style.
$ws.Range("I30:I38").NumberFormat = "@"
$ws.Range("K30:K38").NumberFormat = "@"

# Row 30
$ws.Range("A30").Value = "BTC/USDC"
$ws.Range("B30").Value = 110073.37
$ws.Range("C30").Value = 110264.14734
$ws.Range("D30").Value = 0.001
$ws.Range("E30").Value = 0.1907773400000005
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0.1907773400000005
$ws.Range("H30").Value = 0.1733183421203517
$ws.Range("I30").Value = "2025-06-10"
$ws.Range("J30").Value = "02:16:38"
$ws.Range("K30").Value = "2025-06-10"
$ws.Range("L30").Value = "02:26:57"
$ws.Range("M30").Value = 10.3220522

# Row 31
$ws.Range("A31").Value = "BTC/USDC"
$ws.Range("B31").Value = 109850.71
$ws.Range("C31").Value = 109576.95
$ws.Range("D31").Value = 0.001
$ws.Range("E31").Value = -0.2737600000000093
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = -0.2737600000000093
$ws.Range("H31").Value = -0.2492109518454722
$ws.Range("I31").Value = "2025-06-10"
$ws.Range("J31").Value = "03:57:39"
$ws.Range("K31").Value = "2025-06-10"
$ws.Range("L31").Value = "04:05:57"
$ws.Range("M31").Value = 8.303035916666667

# Row 32
$ws.Range("A32").Value = "BTC/USDC"
$ws.Range("B32").Value = 109606.78947
$ws.Range("C32").Value = 109675.18
$ws.Range("D32").Value = 0.001
$ws.Range("E32").Value = 0.06839052999998967
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0.06839052999998967
$ws.Range("H32").Value = 0.06239625330756408
$ws.Range("I32").Value = "2025-06-10"
$ws.Range("J32").Value = "04:14:39"
$ws.Range("K32").Value = "2025-06-10"
$ws.Range("L32").Value = "04:27:18"
$ws.Range("M32").Value = 12.6649506

# Row 33
$ws.Range("A33").Value = "BTC/USDC"
$ws.Range("B33").Value = 109538.4191
$ws.Range("C33").Value = 109623.4
$ws.Range("D33").Value = 0.001
$ws.Range("E33").Value = 0.08498089999999502
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0.08498089999999502
$ws.Range("H33").Value = 0.07758090786613792
$ws.Range("I33").Value = "2025-06-10"
$ws.Range("J33").Value = "05:20:39"
$ws.Range("K33").Value = "2025-06-10"
$ws.Range("L33").Value = "06:02:48"
$ws.Range("M33").Value = 42.15393601666666

# Row 34
$ws.Range("A34").Value = "BTC/USDC"
$ws.Range("B34").Value = 109421.77
$ws.Range("C34").Value = 109566.47677
$ws.Range("D34").Value = 0.001
$ws.Range("E34").Value = 0.1447067699999898
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0.1447067699999898
$ws.Range("H34").Value = 0.1322467823359006
$ws.Range("I34").Value = "2025-06-10"
$ws.Range("J34").Value = "07:41:39"
$ws.Range("K34").Value = "2025-06-10"
$ws.Range("L34").Value = "08:36:48"
$ws.Range("M34").Value = 55.14536743333333

# Row 35
$ws.Range("A35").Value = "BTC/USDC"
$ws.Range("B35").Value = 109461.97
$ws.Range("C35").Value = 109223.09799
$ws.Range("D35").Value = 0.001
$ws.Range("E35").Value = -0.2388720100000064
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = -0.2388720100000064
$ws.Range("H35").Value = -0.2182237447398457
$ws.Range("I35").Value = "2025-06-10"
$ws.Range("J35").Value = "09:13:38"
$ws.Range("K35").Value = "2025-06-10"
$ws.Range("L35").Value = "09:37:38"
$ws.Range("M35").Value = 23.99959545

# Row 36
$ws.Range("A36").Value = "BTC/USDC"
$ws.Range("B36").Value = 109211.47936
$ws.Range("C36").Value = 109399.35
$ws.Range("D36").Value = 0.001
$ws.Range("E36").Value = 0.1878706400000083
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0.1878706400000083
$ws.Range("H36").Value = 0.1720246269906478
$ws.Range("I36").Value = "2025-06-10"
$ws.Range("J36").Value = "11:07:38"
$ws.Range("K36").Value = "2025-06-10"
$ws.Range("L36").Value = "11:26:08"
$ws.Range("M36").Value = 18.49002208333333

# Row 37
$ws.Range("A37").Value = "BTC/USDC"
$ws.Range("B37").Value = 109408.938
$ws.Range("C37").Value = 109499.7018
$ws.Range("D37").Value = 0.001
$ws.Range("E37").Value = 0.09076380000000063
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0.09076380000000063
$ws.Range("H37").Value = 0.08295830455826254
$ws.Range("I37").Value = "2025-06-10"
$ws.Range("J37").Value = "13:10:16"
$ws.Range("K37").Value = "2025-06-10"
$ws.Range("L37").Value = "13:12:15"
$ws.Range("M37").Value = 1.973762133333334

# Row 38
$ws.Range("A38").Value = "BTC/USDC"
$ws.Range("B38").Value = 109749.57
$ws.Range("C38").Value = 109450.04
$ws.Range("D38").Value = 0.001
$ws.Range("E38").Value = -0.2995300000000134
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = -0.2995300000000134
$ws.Range("H38").Value = -0.2729213426531087
$ws.Range("I38").Value = "2025-06-10"
$ws.Range("J38").Value = "15:27:51"
$ws.Range("K38").Value = "2025-06-10"
$ws.Range("L38").Value = "15:35:08"
$ws.Range("M38").Value = 7.29405835

$ws.Range("I30:I38").Style = "Normal"
$ws.Range("K30:K38").Style = "Normal"

$ws.Range("A1").Select()